# Updated CHE_grids model - 2025-08-14 22:29
#
# The "solar" worksheet contains a lookup table (columns AC:AG, rows 4-28)
# that maps each distribution process (column AC, e.g. "distr_elc_won-CHE_00xx")
# to the grid cell it connects to (column AG, e.g. "CHE_xx"). This change
# re-assigns the grid_cell (column AG) values for rows 4-28 so each row picks
# up a new grid cell designation, per the authoritative data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$ws.Range("AG4").Value  = "CHE_12"
$ws.Range("AG5").Value  = "CHE_10"
$ws.Range("AG6").Value  = "CHE_22"
$ws.Range("AG7").Value  = "CHE_24"
$ws.Range("AG8").Value  = "CHE_8"
$ws.Range("AG9").Value  = "CHE_5"
$ws.Range("AG10").Value = "CHE_11"
$ws.Range("AG11").Value = "CHE_15"
$ws.Range("AG12").Value = "CHE_25"
$ws.Range("AG13").Value = "CHE_3"
$ws.Range("AG14").Value = "CHE_13"
$ws.Range("AG15").Value = "CHE_9"
$ws.Range("AG16").Value = "CHE_21"
$ws.Range("AG17").Value = "CHE_4"
$ws.Range("AG18").Value = "CHE_20"
$ws.Range("AG19").Value = "CHE_1"
$ws.Range("AG20").Value = "CHE_6"
$ws.Range("AG21").Value = "CHE_0"
$ws.Range("AG22").Value = "CHE_7"
$ws.Range("AG23").Value = "CHE_17"
$ws.Range("AG24").Value = "CHE_19"
$ws.Range("AG25").Value = "CHE_23"
$ws.Range("AG26").Value = "CHE_2"
$ws.Range("AG27").Value = "CHE_14"
$ws.Range("AG28").Value = "CHE_18"
